$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of training params (rows 36-42), matching the existing layout:
# A=data_words_type, B=num_topics, C=use_bigram, D=use_trigram,
# E=filter_no_above, F=alpha_entry, G=eta_entry, H=random_state,
# I=passes, J=topn
$rows = @(
    @("qa", 50, $false, $false, 0.2, 0.1, "auto", 1000, 2, 20),
    @("qa", 50, $false, $true,  0.2, 0.1, "auto", 1000, 2, 20),
    @("qa", 50, $true,  $false, 0.3, 0.1, "auto", 1000, 2, 20),
    @("qa", 50, $true,  $false, 0.4, 0.1, "auto", 1000, 2, 20),
    @("qa", 50, $true,  $false, 0.5, 0.1, "auto", 1000, 2, 20),
    @("qa", 50, $true,  $false, 0.6, 0.1, "auto", 1000, 2, 20),
    @("qa", 50, $true,  $false, 0.7, 0.1, "auto", 1000, 2, 20)
)

$startRow = 36
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}

# Mirror the author's final selection / scroll position from the commit.
$ws.Range("E38").Select()
